# Atualização de bases das ligas, do dia: 12-04-2024 às 20:28
# Updates row 129 (previously-unplayed match now has a result + new odds)
# and appends a brand-new row 130 (new upcoming fixture).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 129 : update existing fixture (result now known, odds refreshed)
# ---------------------------------------------------------------------
$ws.Range("B129").Value = 7749763
$ws.Range("E129").Value = 45393.45833333334
$ws.Range("F129").Value = "Bengaluru"
$ws.Range("G129").Value = "Mohun Bagan SG"
$ws.Range("H129").Value = 0
$ws.Range("I129").Value = 4
$ws.Range("J129").Value = "A"
$ws.Range("K129").Value = 4.333
$ws.Range("L129").Value = 3.6
$ws.Range("M129").Value = 1.8
$ws.Range("N129").Value = 3.8
$ws.Range("O129").Value = 3.6
$ws.Range("P129").Value = 1.85
$ws.Range("Q129").Value = 0.5
$ws.Range("R129").Value = 1.95
$ws.Range("S129").Value = 1.9
$ws.Range("T129").Value = 2.75
$ws.Range("U129").Value = 2
$ws.Range("V129").Value = 1.85
$ws.Range("W129").Value = -1
$ws.Range("X129").Value = -1
$ws.Range("Y129").Value = 0.8500000000000001
$ws.Range("Z129").Value = -1
$ws.Range("AA129").Value = 0.8999999999999999
$ws.Range("AB129").Value = 1
$ws.Range("AC129").Value = -1

# ---------------------------------------------------------------------
# Row 130 : brand-new fixture row
# ---------------------------------------------------------------------
$ws.Range("A130").Value = 128
$ws.Range("B130").Value = 7749775
$ws.Range("C130").Value = "India Super League"
$ws.Range("D130").Value = "India Super League"
$ws.Range("E130").Value = 45395.45833333334
$ws.Range("F130").Value = "Northeast United"
$ws.Range("G130").Value = "Odisha FC"
$ws.Range("K130").Value = 2.8
$ws.Range("L130").Value = 3.5
$ws.Range("M130").Value = 2.2
$ws.Range("N130").Value = 2.9
$ws.Range("O130").Value = 3.5
$ws.Range("P130").Value = 2.1
$ws.Range("Q130").Value = 0.25
$ws.Range("R130").Value = 1.9
$ws.Range("S130").Value = 1.9
$ws.Range("T130").Value = 2.75
$ws.Range("U130").Value = 1.8
$ws.Range("V130").Value = 2
$ws.Range("W130").Value = 0
$ws.Range("X130").Value = 0
$ws.Range("Y130").Value = 0
$ws.Range("Z130").Value = 0
$ws.Range("AA130").Value = 0

# Match the formatting of the row above (bold/bordered id column, date format)
$ws.Range("A129").Copy()
$ws.Range("A130").PasteSpecial(-4122)
$ws.Range("E129").Copy()
$ws.Range("E130").PasteSpecial(-4122)
$excel.CutCopyMode = $false
